# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.278.06"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "'1.666.14"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("D5").Value = "'218.44"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.5336"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'0.2639"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "'0.06372"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'20.55"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "'0.07820"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'4.570"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "'1.670.27"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'1.893.75"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'0.5538"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "'0.0₅8208"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'65.80"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'4.689"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'194.20"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'10.22"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'6.041"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'1.010"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'145.96"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").Value = "'0.1229"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'16.16"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("D29").Value = "'0.05880"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'1.281"
$ws.Range("D31").Value = "'3.596"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").Value = "'3.279"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "'1.610"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").Value = "'0.9614"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'2.823"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "'0.5797"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").Value = "'0.01607"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").Value = "'0.8633"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "'5.847"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.049.47"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.009"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "'104.00"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "'1.803.53"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'57.65"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("D48").Value = "'0.4377"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "'8.041"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("D50").Value = "'0.05159"
$ws.Range("E51").Value = "  -3.28%  "
